# Implemented POM design pattern in framework
# - TC_UI_001: locators for the "addUser" flow now reference page-object
#   members (DashboardPage.addUser / DashboardPage.name / DashboardPage.save)
#   instead of raw CSS ids.
# - TC_UI_002: the assertion locator now references the page-object member
#   ExamplePage.title instead of the raw "h1" selector.
# - TC_UI_001 becomes the active/selected worksheet when the workbook is saved.

$wb = $excel.ActiveWorkbook

$wsTcUi001 = $wb.Worksheets.Item("TC_UI_001")
$wsTcUi002 = $wb.Worksheets.Item("TC_UI_002")

# --- TC_UI_001: switch raw CSS locators to Page Object Model references ---
$wsTcUi001.Range("C3").Value = "DashboardPage.addUser"
$wsTcUi001.Range("C4").Value = "DashboardPage.name"
$wsTcUi001.Range("C5").Value = "DashboardPage.save"

# Column C now needs to be wide enough to show the longer POM locator text.
$wsTcUi001.Range("C:C").ColumnWidth = 30.7109375

# --- TC_UI_002: switch the assertion locator to its Page Object reference ---
$wsTcUi002.Range("C3").Value = "ExamplePage.title"

# Leave the cursor where the author left it after editing the flow.
$wsTcUi001.Range("C9").Select()

# Make TC_UI_001 the active sheet/tab, since that's where the POM work happened.
$wsTcUi001.Activate()
